$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.02025219827703495
$ws.Range("D2").Value = 0.009188019191155306
$ws.Range("E2").Value = 0.06840950260453127
$ws.Range("F2").Value = 4.535714456049959
$ws.Range("G2").Value = 0.002622642932739561
$ws.Range("I2").Value = 1.275920412856571
$ws.Range("J2").Value = 0.1750878153708513
$ws.Range("K2").Value = 2.024989549009234
$ws.Range("M2").Value = 0.5335331956844556
$ws.Range("B3").Value = 0.01788570639442355
$ws.Range("D3").Value = 0.009125331924932922
$ws.Range("E3").Value = 0.06855475687318702
$ws.Range("F3").Value = 4.471259589366497
$ws.Range("G3").Value = 0.002628467586084825
$ws.Range("I3").Value = 1.261105888783113
$ws.Range("J3").Value = 0.1741598512296605
$ws.Range("K3").Value = 1.949217145640148
$ws.Range("M3").Value = 0.5207405224949468
$ws.Range("B4").Value = 0.01642031159172319
$ws.Range("D4").Value = 0.00909139403914061
$ws.Range("E4").Value = 0.06867139211482787
$ws.Range("F4").Value = 4.433416738949489
$ws.Range("G4").Value = 0.002632230520959587
$ws.Range("I4").Value = 1.252247528314491
$ws.Range("J4").Value = 0.1736507215066538
$ws.Range("K4").Value = 1.904203304309135
$ws.Range("M4").Value = 0.5132599846696024
$ws.Range("B5").Value = 0.015820145178985
$ws.Range("D5").Value = 0.00907870342546957
$ws.Range("E5").Value = 0.06872583519260722
$ws.Range("F5").Value = 4.418429246569502
$ws.Range("G5").Value = 0.002633811029218306
$ws.Range("I5").Value = 1.248697003065168
$ws.Range("J5").Value = 0.1734584699738591
$ws.Range("K5").Value = 1.886237921082028
$ws.Range("M5").Value = 0.5103054425265299
$ws.Range("B6").Value = 0.01572030943857783
$ws.Range("D6").Value = 0.009076664817731839
$ws.Range("E6").Value = 0.0687352933006693
$ws.Range("F6").Value = 4.415966728131508
$ws.Range("G6").Value = 0.002634076319891192
$ws.Range("I6").Value = 1.248111010543646
$ws.Range("J6").Value = 0.1734274654245951
$ws.Range("K6").Value = 1.88327756995497
$ws.Range("M6").Value = 0.5098205031493137
$ws.Range("B7").Value = 0.01641222957409383
$ws.Range("D7").Value = 0.009091218283057145
$ws.Range("E7").Value = 0.06867209834646815
$ws.Range("F7").Value = 4.433212858619811
$ws.Range("G7").Value = 0.002632251645332929
$ws.Range("I7").Value = 1.252199405106246
$ws.Range("J7").Value = 0.1736480671319924
$ws.Range("K7").Value = 1.90395948817806
$ws.Range("M7").Value = 0.5132197590790142
$ws.Range("B8").Value = 0.01943885717162885
$ws.Range("D8").Value = 0.009165456482897127
$ws.Range("E8").Value = 0.06845389581608075
$ws.Range("F8").Value = 4.513129604977081
$ws.Range("G8").Value = 0.002624612650010074
$ws.Range("I8").Value = 1.270762716520863
$ws.Range("J8").Value = 0.1747552499858003
$ws.Range("K8").Value = 1.998548573759024
$ws.Range("M8").Value = 0.5290444480914118
$ws.Range("B9").Value = 0.02527179798854462
$ws.Range("D9").Value = 0.00934742102553443
$ws.Range("E9").Value = 0.06824336557326305
$ws.Range("F9").Value = 4.683693395193899
$ws.Range("G9").Value = 0.002611105391631944
$ws.Range("I9").Value = 1.30907571751986
$ws.Range("J9").Value = 0.1774090914594808
$ws.Range("K9").Value = 2.196117998033571
$ws.Range("M9").Value = 0.5630604423179051
$ws.Range("B10").Value = 0.02948965397224868
$ws.Range("D10").Value = 0.009503684339815521
$ws.Range("E10").Value = 0.06822071569374621
$ws.Range("F10").Value = 4.817603314499962
$ws.Range("G10").Value = 0.002602068809974685
$ws.Range("I10").Value = 1.338424411981165
$ws.Range("J10").Value = 0.1796555645133964
$ws.Range("K10").Value = 2.348787077433144
$ws.Range("M10").Value = 0.5898952716009092
$ws.Range("B11").Value = 0.03139269984411186
$ws.Range("D11").Value = 0.009579762727199181
$ws.Range("E11").Value = 0.06823898856097621
$ws.Range("F11").Value = 4.880424045692934
$ws.Range("G11").Value = 0.00259814820486197
$ws.Range("I11").Value = 1.352044329464192
$ws.Range("J11").Value = 0.1807425348816878
$ws.Range("K11").Value = 2.419906352706334
$ws.Range("M11").Value = 0.6025090107102926
$ws.Range("B12").Value = 0.03211098132236145
$ws.Range("D12").Value = 0.009609296571255754
$ws.Range("E12").Value = 0.06825000894782285
$ws.Range("F12").Value = 4.904489014340129
$ws.Range("G12").Value = 0.002596690747246045
$ws.Range("I12").Value = 1.357241096247435
$ws.Range("J12").Value = 0.1811635346093396
$ws.Range("K12").Value = 2.447079904043733
$ws.Range("M12").Value = 0.6073443498402682
$ws.Range("B13").Value = 0.03195639323151056
$ws.Range("D13").Value = 0.009602903595755663
$ws.Range("E13").Value = 0.06824745325452319
$ws.Range("F13").Value = 4.899293875726528
$ws.Range("G13").Value = 0.002597003429952764
$ws.Range("I13").Value = 1.356120127965994
$ws.Range("J13").Value = 0.1810724466980673
$ws.Range("K13").Value = 2.441216785290919
$ws.Range("M13").Value = 0.6063003527494004
$ws.Range("B14").Value = 0.03145184121153477
$ws.Range("D14").Value = 0.009582177933161518
$ws.Range("E14").Value = 0.06823981307469218
$ws.Range("F14").Value = 4.882398338164194
$ws.Range("G14").Value = 0.002598027754769722
$ws.Range("I14").Value = 1.352471082043735
$ws.Range("J14").Value = 0.1807769824356811
$ws.Range("K14").Value = 2.422137069913845
$ws.Range("M14").Value = 0.6029056367233991
$ws.Range("B15").Value = 0.0311424779568199
$ws.Range("D15").Value = 0.009569577429228815
$ws.Range("E15").Value = 0.06823566705157624
$ws.Range("F15").Value = 4.872085360532395
$ws.Range("G15").Value = 0.002598658719912097
$ws.Range("I15").Value = 1.350241056193056
$ws.Range("J15").Value = 0.1805972256512405
$ws.Range("K15").Value = 2.410481803730192
$ws.Range("M15").Value = 0.600833941716985
$ws.Range("B16").Value = 0.02936496087475149
$ws.Range("D16").Value = 0.009498813476788825
$ws.Range("E16").Value = 0.06822009582369049
$ws.Range("F16").Value = 4.813536346942158
$ws.Range("G16").Value = 0.002602328844057601
$ws.Range("I16").Value = 1.33753977590419
$ws.Range("J16").Value = 0.179585840171633
$ws.Range("K16").Value = 2.344173041728652
$ws.Range("M16").Value = 0.5890791418870194
$ws.Range("B17").Value = 0.02827042044670947
$ws.Range("D17").Value = 0.009456685647450058
$ws.Range("E17").Value = 0.06821785640275202
$ws.Range("F17").Value = 4.778107795947307
$ws.Range("G17").Value = 0.002604628944124832
$ws.Range("I17").Value = 1.329817224463582
$ws.Range("J17").Value = 0.1789820706117169
$ws.Range("K17").Value = 2.303923853679009
$ws.Range("M17").Value = 0.5819722936157063
$ws.Range("B18").Value = 0.0276394003678746
$ws.Range("D18").Value = 0.009432924437454915
$ws.Range("E18").Value = 0.06821925857709488
$ws.Range("F18").Value = 4.75790934751106
$ws.Range("G18").Value = 0.002605969810853806
$ws.Range("I18").Value = 1.325400710663274
$ws.Range("J18").Value = 0.1786409179057955
$ws.Range("K18").Value = 2.280930578491166
$ws.Range("M18").Value = 0.5779228573386845
$ws.Range("B19").Value = 0.02742549869137889
$ws.Range("D19").Value = 0.009424959758536744
$ws.Range("E19").Value = 0.06822019564011761
$ws.Range("F19").Value = 4.751101194990952
$ws.Range("G19").Value = 0.002606426886287908
$ws.Range("I19").Value = 1.323909683635051
$ws.Range("J19").Value = 0.1785264594465019
$ws.Range("K19").Value = 2.27317234893826
$ws.Range("M19").Value = 0.5765583443574087
$ws.Range("B20").Value = 0.02838708898993758
$ws.Range("D20").Value = 0.009461121577153619
$ws.Range("E20").Value = 0.06821781642617708
$ws.Range("F20").Value = 4.781860671458929
$ws.Range("G20").Value = 0.002604382242133938
$ws.Range("I20").Value = 1.330636681771622
$ws.Range("J20").Value = 0.179045709307097
$ws.Range("K20").Value = 2.308192184676386
$ws.Range("M20").Value = 0.5827248701236911
$ws.Range("B21").Value = 0.03160010528254276
$ws.Range("D21").Value = 0.009588245844877719
$ws.Range("E21").Value = 0.06824194595022703
$ws.Range("F21").Value = 4.887353455013482
$ws.Range("G21").Value = 0.002597726149119696
$ws.Range("I21").Value = 1.353541827337253
$ws.Range("J21").Value = 0.1808635123891094
$ws.Range("K21").Value = 2.427734653279117
$ws.Range("M21").Value = 0.6039011487713637
$ws.Range("B22").Value = 0.03368618836637438
$ws.Range("D22").Value = 0.009675555631222466
$ws.Range("E22").Value = 0.06828161520226317
$ws.Range("F22").Value = 4.957909610265204
$ws.Range("G22").Value = 0.002593534421064751
$ws.Range("I22").Value = 1.368740388093144
$ws.Range("J22").Value = 0.1821062871064711
$ws.Range("K22").Value = 2.507275490975132
$ws.Range("M22").Value = 0.6180839110104799
$ws.Range("B23").Value = 0.03257410501529279
$ws.Range("D23").Value = 0.009628567805521016
$ws.Range("E23").Value = 0.06825825888934212
$ws.Range("F23").Value = 4.920104341273202
$ws.Range("G23").Value = 0.002595757182350722
$ws.Range("I23").Value = 1.360607532585888
$ws.Range("J23").Value = 0.1814379747057089
$ws.Range("K23").Value = 2.464693032895241
$ws.Range("M23").Value = 0.6104828174194665
$ws.Range("B24").Value = 0.02833434859756778
$ws.Range("D24").Value = 0.00945911466489946
$ws.Range("E24").Value = 0.06821782612130534
$ws.Range("F24").Value = 4.780163467417964
$ws.Range("G24").Value = 0.002604493718380207
$ws.Range("I24").Value = 1.330266132637092
$ws.Range("J24").Value = 0.1790169196837681
$ws.Range("K24").Value = 2.306262016021776
$ws.Range("M24").Value = 0.5823845169776405
$ws.Range("B25").Value = 0.02370532763469413
$ws.Range("D25").Value = 0.009294259448711983
$ws.Range("E25").Value = 0.06827710045995516
$ws.Range("F25").Value = 4.636053019398958
$ws.Range("G25").Value = 0.002614602894250567
$ws.Range("I25").Value = 1.29850285000515
$ws.Range("J25").Value = 0.1766392344056769
$ws.Range("K25").Value = 2.141361595000831
$ws.Range("M25").Value = 0.5535363054245224
